$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 3) mirroring the structure of row 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "cent/Wh"

# Auto-fit the columns so widths match the new content (as seen in the diff)
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 17.833333333333336

# Select the last edited cell, matching the recorded selection in the diff
$ws.Range("E3").Select() | Out-Null
